# Apply the updated values for columns A, B, Q, R in rows 2-8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; A = 112241889; B = 78699; Q = 553253; R = 7007769 },
    @{ Row = 3; A = 112241887; B = 78699; Q = 553188; R = 7007666 },
    @{ Row = 4; A = 112241888; B = 78699; Q = 553253; R = 7007768 },
    @{ Row = 5; A = 112241875; B = 78726; Q = 553188; R = 7007668 },
    @{ Row = 6; A = 112241886; B = 78699; Q = 553306; R = 7007600 },
    @{ Row = 7; A = 112241884; B = 78699; Q = 553337; R = 7007616 },
    @{ Row = 8; A = 112241885; B = 78699; Q = 553321; R = 7007611 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("A$r").Value = $u.A
    $ws.Range("B$r").Value = $u.B
    $ws.Range("Q$r").Value = $u.Q
    $ws.Range("R$r").Value = $u.R
}
